$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.356.16'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '3.968.16'
$ws.Range('E3').Value = '  -2.00%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Formula = "'612.51"
$ws.Range('E5').Value = '  +8.28%  '
$ws.Range('D6').Formula = "'169.03"
$ws.Range('E6').Value = '  +11.73%  '
$ws.Range('D7').Formula = "'0.679"
$ws.Range('E7').Value = '  -2.67%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').Formula = "'0.776"
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').Formula = "'0.186"
$ws.Range('E10').Value = '  +7.88%  '
$ws.Range('D11').Formula = "'56.23"
$ws.Range('E11').Value = '  +3.81%  '
$ws.Range('D12').Formula = "'0.0000337"
$ws.Range('E12').Value = '  +1.75%  '
$ws.Range('D13').Formula = "'11.26"
$ws.Range('E13').Value = '  +1.73%  '
$ws.Range('D14').Value = '4.603.36'
$ws.Range('E14').Value = '  -1.90%  '
$ws.Range('D15').Value = '3.977.23'
$ws.Range('E15').Value = '  -1.72%  '
$ws.Range('D16').Formula = "'14.17"
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('D17').Formula = "'1.24"
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('D18').Formula = "'20.48"
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').Value = '73.186.63'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').Formula = "'440.89"
$ws.Range('E21').Value = '  -1.40%  '
$ws.Range('D22').Formula = "'4.84"
$ws.Range('E22').Value = '  +8.01%  '
$ws.Range('D23').Formula = "'95.32"
$ws.Range('E23').Value = '  -3.03%  '
$ws.Range('D24').Formula = "'3.44"
$ws.Range('E24').Value = '  -2.98%  '
$ws.Range('D25').Formula = "'14.20"
$ws.Range('E25').Value = '  -4.02%  '
$ws.Range('D26').Formula = "'4.19"
$ws.Range('E26').Value = '  -4.29%  '
$ws.Range('E27').Value = '  -1.85%  '
$ws.Range('D28').Formula = "'5.96"
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').Formula = "'10.50"
$ws.Range('E29').Value = '  -4.56%  '
$ws.Range('D30').Formula = "'36.04"
$ws.Range('E30').Value = '  -3.19%  '
$ws.Range('D31').Formula = "'7.89"
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').Formula = "'13.88"
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('D33').Formula = "'0.0000106"
$ws.Range('E33').Value = '  +14.21%  '
$ws.Range('E34').Value = '  -4.51%  '
$ws.Range('D35').Formula = "'48.16"
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('D36').Formula = "'70.45"
$ws.Range('E36').Value = '  +3.95%  '
$ws.Range('D37').Formula = "'649.10"
$ws.Range('E37').Value = '  -5.82%  '
$ws.Range('D38').Formula = "'0.431"
$ws.Range('E38').Value = '  -3.93%  '
$ws.Range('D39').Formula = "'3.40"
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('D40').Formula = "'1.00"
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('D41').Formula = "'0.146"
$ws.Range('E41').Value = '  -2.54%  '
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').Formula = "'0.0483"
$ws.Range('E43').Value = '  -2.70%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Formula = "'3.17"
$ws.Range('E44').Value = '  +35.84%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Formula = "'3.16"
$ws.Range('E45').Value = '  -5.95%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').Formula = "'10.44"
$ws.Range('E46').Value = '  -5.85%  '
$ws.Range('E47').Value = '  -2.33%  '
$ws.Range('D48').Formula = "'0.000303"
$ws.Range('E48').Value = '  +6.91%  '
$ws.Range('D49').Formula = "'3.41"
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('D50').Formula = "'2.57"
$ws.Range('E50').Value = '  -4.20%  '
$ws.Range('D51').Formula = "'3.00"
$ws.Range('E51').Value = '  -3.88%  '
